# "Made a change to the schema and subsequent changes to the server code"
#
# The DB schema sheet (Sheet1) has a "Cholsterol" column (column I, row 4)
# that is being removed from the schema entirely. Select the whole column
# and delete it, which shifts every column to its right (Sodium, Total
# Carbohydrate, Dietary Fiber, Sugars, Protein) one position to the left
# and drops the now-unused "Cholsterol" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns("I").Select()
$ws.Columns("I").Delete()
